$d = $word.ActiveDocument

# 1. "factura" table: the nif column loses its "unique" constraint, keeping
#    only "not null". (wdFindContinue=1, wdReplaceAll=2)
$r = $d.Content
$null = $r.Find.Execute("nif char(9) unique not null,", $true, $false, $false, `
    $false, $false, $true, 1, $false, "nif char(9) not null,", 2)

# 2. Re-locate the edited text (now the only "nif char(9) not null," in the
#    document) and split the run right after "nif char(9) " so a collapsed
#    "_GoBack" bookmark can sit exactly between the two pieces - this is
#    Word's usual behaviour of re-stamping _GoBack at the last edit point.
$prefix = "nif char(9) "
$r2 = $d.Content
$null = $r2.Find.Execute("$prefix" + "not null,", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)

$splitPos = $r2.Start + $prefix.Length
$bmRange = $d.Range($splitPos, $splitPos)

# Adding a bookmark with a name that already exists elsewhere in the document
# moves it, so this single call both removes the old "_GoBack" (that used to
# sit next to "Freddy Fazbear") and creates the new one here.
$d.Bookmarks.Add("_GoBack", $bmRange)
